# "Mas mediciones de temperatura"
# Replace the raw-data block on "Datos crudos" with a new, longer run (38
# samples instead of 35) and let the dependent formulas on "Datos validos"
# recompute from it. Also flips which sheet/cell is selected/active and
# nudges the two +/-0.5 offset formulas (were +/-0.6) and the chart-1 value
# axis minimum (17 -> 19).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Datos crudos")
$ws2 = $wb.Worksheets.Item("Datos válidos")

# ---------------------------------------------------------------------
# 1) New timestamp (col C) / temperature (col E) series, rows 2..39.
#    Rows 2-36 overwrite the previous readings in place; rows 37-39 are
#    brand new rows appended at the end of the raw-data table.
# ---------------------------------------------------------------------
$timestamps = @(
    "2023-12-09 06:23:18","2023-12-09 06:24:20","2023-12-09 06:25:22","2023-12-09 06:26:24",
    "2023-12-09 06:27:26","2023-12-09 06:28:28","2023-12-09 06:29:30","2023-12-09 06:30:32",
    "2023-12-09 06:31:35","2023-12-09 06:32:37","2023-12-09 06:33:39","2023-12-09 06:34:41",
    "2023-12-09 06:35:43","2023-12-09 06:36:45","2023-12-09 06:37:47","2023-12-09 06:38:50",
    "2023-12-09 06:39:52","2023-12-09 06:40:54","2023-12-09 06:41:56","2023-12-09 06:42:58",
    "2023-12-09 06:44:00","2023-12-09 06:45:02","2023-12-09 06:46:04","2023-12-09 06:47:07",
    "2023-12-09 06:48:09","2023-12-09 06:49:11","2023-12-09 06:50:13","2023-12-09 06:51:15",
    "2023-12-09 06:52:17","2023-12-09 06:53:19","2023-12-09 06:54:22","2023-12-09 06:55:24",
    "2023-12-09 06:56:26","2023-12-09 06:57:28","2023-12-09 06:58:30","2023-12-09 06:59:32",
    "2023-12-09 07:00:34","2023-12-09 07:01:36"
)

$temperatures = @(
    26.8531468531468,  25.673076923076898, 24.7552447552447,   23.181818181818102,
    21.870629370629299,21.477272727272702, 21.215034965034899, 21.215034965034899,
    20.428321678321701,20.428321678321701, 20.2972027972028,   20.166083916083899,
    20.2972027972028,  20.2972027972028,    20.166083916083899, 20.034965034965001,
    20.034965034965001,19.9038461538461,   20.428321678321701, 20.034965034965001,
    19.9038461538461,  19.772727272727298,  19.772727272727298, 20.034965034965001,
    19.772727272727298,19.9038461538461,    19.510489510489499, 19.9038461538461,
    19.9038461538461,  19.510489510489499,  19.772727272727298, 19.641608391608401,
    19.641608391608401,19.641608391608401,  19.641608391608401, 20.034965034965001,
    19.641608391608401,19.510489510489499
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2

    if ($row -gt 36) {
        # Brand-new row: fill in the A/B/D boilerplate columns too, and force
        # the timestamp cell to the same Text number format used by the rest
        # of column C (style index 5 in the original file).
        $ws1.Range("A$row").Value = 23
        $ws1.Range("B$row").Value = 17
        $ws1.Range("C$row").NumberFormat = "@"
        $ws1.Range("D$row").Value = 0
    }

    $ws1.Range("C$row").Value = $timestamps[$i]
    $ws1.Range("E$row").Value = $temperatures[$i]
}

# ---------------------------------------------------------------------
# 2) Housekeeping formulas on "Datos crudos" that key off the new extent.
# ---------------------------------------------------------------------
$ws1.Range("H2").Value = 11
$ws1.Range("H3").Formula = "=COUNT(E:E)"

# ---------------------------------------------------------------------
# 3) The two offset formulas on "Datos validos" (+/-0.6 -> +/-0.5).
# ---------------------------------------------------------------------
$ws2.Range("E17").Formula = "=E14+0.5"
$ws2.Range("E20").Formula = "=E14-0.5"

# ---------------------------------------------------------------------
# 4) Chart 1 (on "Datos crudos") value-axis minimum: 17 -> 19.
# ---------------------------------------------------------------------
$chart1 = $ws1.ChartObjects().Item(1).Chart
$valueAxis1 = $chart1.Axes(2, 1)
$valueAxis1.MinimumScale = 19

# ---------------------------------------------------------------------
# 5) Selection / active-sheet swap: "Datos crudos" becomes the active tab
#    (selection H3) and "Datos validos" loses tabSelected (selection B3).
# ---------------------------------------------------------------------
$ws2.Range("B3").Select()
$ws1.Activate()
$ws1.Range("H3").Select()
